$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview": the two data rows (6e69aeeb... and 91585311...)
# swap places, and the 6e69aeeb row's status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" with an
# updated "Latest HO Xliff Generate Date" timestamp.
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md"
$wsOverview.Range("A3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md", "", "", "e2e\91585311-1d2e-4f85-8298-44c6a5eb3f6d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/91585311-1d2e-4f85-8298-44c6a5eb3f6d.md", "", "", "e2e\6e69aeeb-2812-4203-93f3-f3326070ca5b.md")

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 13:01:04"

# -----------------------------------------------------------------
# Sheet "zh-cn": row 2 becomes the 91585311 file (still handed back,
# in sync), row 3 becomes the 6e69aeeb file which is now ready for a
# new handoff (status/error/timestamps updated).
# -----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md"
$wsZhCn.Range("G2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.91ff2c4f4910992f58f4d42937a1cb1a4402e2fa.zh-cn.xlf"
$wsZhCn.Range("I2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md"
$wsZhCn.Range("J2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.91ff2c4f4910992f58f4d42937a1cb1a4402e2fa.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.1f93288f09a41b6a7e9fdf63ebb8edaaec3e2ee8.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-13 13:00:54"
$wsZhCn.Range("I3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.md"
$wsZhCn.Range("J3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.1f93288f09a41b6a7e9fdf63ebb8edaaec3e2ee8.zh-cn.xlf"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e08fae468278915ec5f2a192294f6dc9ddccfee0/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md."

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Range("I3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md", "", "", "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7402f124d08003221017d7ca3ccd46b5657b1c87/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md", "", "", "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/91585311-1d2e-4f85-8298-44c6a5eb3f6d.md", "", "", "6e69aeeb-2812-4203-93f3-f3326070ca5b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7402f124d08003221017d7ca3ccd46b5657b1c87/e2e/91585311-1d2e-4f85-8298-44c6a5eb3f6d.md", "", "", "6e69aeeb-2812-4203-93f3-f3326070ca5b.md")

$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# -----------------------------------------------------------------
# Sheet "de-de": same pattern as zh-cn above, with the de-de xliff
# file names / timestamps.
# -----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md"
$wsDeDe.Range("G2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.91ff2c4f4910992f58f4d42937a1cb1a4402e2fa.de-de.xlf"
$wsDeDe.Range("I2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md"
$wsDeDe.Range("J2").Value = "91585311-1d2e-4f85-8298-44c6a5eb3f6d.91ff2c4f4910992f58f4d42937a1cb1a4402e2fa.de-de.xlf"

$wsDeDe.Range("A3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.1f93288f09a41b6a7e9fdf63ebb8edaaec3e2ee8.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-13 13:01:04"
$wsDeDe.Range("I3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.md"
$wsDeDe.Range("J3").Value = "6e69aeeb-2812-4203-93f3-f3326070ca5b.1f93288f09a41b6a7e9fdf63ebb8edaaec3e2ee8.de-de.xlf"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e08fae468278915ec5f2a192294f6dc9ddccfee0/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md."

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Range("I3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md", "", "", "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d816d1a56a400f8ec2d4400fc6586606763ae256/e2e/6e69aeeb-2812-4203-93f3-f3326070ca5b.md", "", "", "91585311-1d2e-4f85-8298-44c6a5eb3f6d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c096d162a72f33d676b6e4af4b648e6e502911b/e2e/91585311-1d2e-4f85-8298-44c6a5eb3f6d.md", "", "", "6e69aeeb-2812-4203-93f3-f3326070ca5b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d816d1a56a400f8ec2d4400fc6586606763ae256/e2e/91585311-1d2e-4f85-8298-44c6a5eb3f6d.md", "", "", "6e69aeeb-2812-4203-93f3-f3326070ca5b.md")

$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
